# Daily update for 24/05/2020
# Adds the new day's row (24/05/2020, serial 43975) to the three daily
# data tables and refreshes the "last active cell" selection on the
# sheets whose tables grew.

$wb = $excel.ActiveWorkbook

$sheetContents = $wb.Worksheets.Item(1)
$sheetCases    = $wb.Worksheets.Item(3)   # Table 1 - Cumulative cases
$sheetICU      = $wb.Worksheets.Item(4)   # Table 2 - ICU patients
$sheetConf     = $wb.Worksheets.Item(5)   # Table 3a - Hospital Confirmed
$sheetSusp     = $wb.Worksheets.Item(6)   # Table 3b - Hospital Suspected

# --- style "donor" cells: cells elsewhere in the workbook that already
# carry the exact cellXf we need for the new row, so copying their
# format (rather than poking Font/Interior/NumberFormat by hand) is
# guaranteed to land on the same style index instead of synthesising a
# near-duplicate one. ---
$donor5  = $sheetCases.Range("A5")     # date, header-row style
$donor28 = $sheetCases.Range("B47")    # plain number cell
$donor41 = $sheetCases.Range("O58")    # number cell w/ extra border (last col before total)
$donor32 = $sheetCases.Range("P5")     # total column number cell
$donor30 = $sheetICU.Range("A36")      # date cell (tables 2/3a/3b style)
$donor14 = $sheetICU.Range("B4")       # "*"/text cell
$donor50 = $sheetICU.Range("B57")      # number cell, no-fill variant
$donor37 = $sheetICU.Range("P37")      # "*" cell in the total column
$donor13 = $sheetICU.Range("Q56")      # total column number cell

function Set-CellFromDonor($cell, $donor, $value) {
    $cell.Value2 = $value
    $donor.Copy()
    $cell.PasteSpecial(-4122)   # xlPasteFormats
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Table 1 - Cumulative cases: row 82 (the row already existed, blank)
# ---------------------------------------------------------------------
$r = 82
Set-CellFromDonor $sheetCases.Range("A$r") $donor5  43975
Set-CellFromDonor $sheetCases.Range("B$r") $donor28 1019
Set-CellFromDonor $sheetCases.Range("C$r") $donor28 322
Set-CellFromDonor $sheetCases.Range("D$r") $donor28 258
Set-CellFromDonor $sheetCases.Range("E$r") $donor28 826
Set-CellFromDonor $sheetCases.Range("F$r") $donor28 917
Set-CellFromDonor $sheetCases.Range("G$r") $donor28 1224
Set-CellFromDonor $sheetCases.Range("H$r") $donor28 3876
Set-CellFromDonor $sheetCases.Range("I$r") $donor28 335
Set-CellFromDonor $sheetCases.Range("J$r") $donor28 1934
Set-CellFromDonor $sheetCases.Range("K$r") $donor28 2665
Set-CellFromDonor $sheetCases.Range("L$r") $donor28 7
Set-CellFromDonor $sheetCases.Range("M$r") $donor28 54
Set-CellFromDonor $sheetCases.Range("N$r") $donor28 1658
Set-CellFromDonor $sheetCases.Range("O$r") $donor41 6
Set-CellFromDonor $sheetCases.Range("P$r") $donor32 15101

# ---------------------------------------------------------------------
# Table 2 - ICU patients: new row 71
# ---------------------------------------------------------------------
$r = 71
Set-CellFromDonor $sheetICU.Range("A$r") $donor30 43975
Set-CellFromDonor $sheetICU.Range("B$r") $donor14 "*"
Set-CellFromDonor $sheetICU.Range("C$r") $donor14 "*"
Set-CellFromDonor $sheetICU.Range("D$r") $donor14 "*"
Set-CellFromDonor $sheetICU.Range("E$r") $donor14 "*"
Set-CellFromDonor $sheetICU.Range("F$r") $donor14 "*"
Set-CellFromDonor $sheetICU.Range("G$r") $donor14 "*"
Set-CellFromDonor $sheetICU.Range("H$r") $donor50 12
Set-CellFromDonor $sheetICU.Range("I$r") $donor14 "*"
Set-CellFromDonor $sheetICU.Range("J$r") $donor14 "*"
Set-CellFromDonor $sheetICU.Range("K$r") $donor50 8
Set-CellFromDonor $sheetICU.Range("L$r") $donor14 "*"
Set-CellFromDonor $sheetICU.Range("M$r") $donor14 "*"
Set-CellFromDonor $sheetICU.Range("N$r") $donor14 "*"
Set-CellFromDonor $sheetICU.Range("O$r") $donor14 "*"
Set-CellFromDonor $sheetICU.Range("P$r") $donor37 "*"
Set-CellFromDonor $sheetICU.Range("Q$r") $donor13 44

# ---------------------------------------------------------------------
# Table 3a - Hospital Confirmed: new row 63
# ---------------------------------------------------------------------
$r = 63
Set-CellFromDonor $sheetConf.Range("A$r") $donor30 43975
Set-CellFromDonor $sheetConf.Range("B$r") $donor50 25
Set-CellFromDonor $sheetConf.Range("C$r") $donor50 11
Set-CellFromDonor $sheetConf.Range("D$r") $donor14 "*"
Set-CellFromDonor $sheetConf.Range("E$r") $donor50 69
Set-CellFromDonor $sheetConf.Range("F$r") $donor50 14
Set-CellFromDonor $sheetConf.Range("G$r") $donor50 72
Set-CellFromDonor $sheetConf.Range("H$r") $donor50 379
Set-CellFromDonor $sheetConf.Range("I$r") $donor50 10
Set-CellFromDonor $sheetConf.Range("J$r") $donor50 83
Set-CellFromDonor $sheetConf.Range("K$r") $donor50 165
Set-CellFromDonor $sheetConf.Range("L$r") $donor14 "*"
Set-CellFromDonor $sheetConf.Range("M$r") $donor14 "*"
Set-CellFromDonor $sheetConf.Range("N$r") $donor50 13
Set-CellFromDonor $sheetConf.Range("O$r") $donor14 "*"
Set-CellFromDonor $sheetConf.Range("P$r") $donor37 "*"
Set-CellFromDonor $sheetConf.Range("Q$r") $donor13 845

# ---------------------------------------------------------------------
# Table 3b - Hospital Suspected: new row 63
# ---------------------------------------------------------------------
$r = 63
Set-CellFromDonor $sheetSusp.Range("A$r") $donor30 43975
Set-CellFromDonor $sheetSusp.Range("B$r") $donor50 30
Set-CellFromDonor $sheetSusp.Range("C$r") $donor50 25
Set-CellFromDonor $sheetSusp.Range("D$r") $donor50 18
Set-CellFromDonor $sheetSusp.Range("E$r") $donor50 33
Set-CellFromDonor $sheetSusp.Range("F$r") $donor50 55
Set-CellFromDonor $sheetSusp.Range("G$r") $donor50 18
Set-CellFromDonor $sheetSusp.Range("H$r") $donor14 "N/A"
Set-CellFromDonor $sheetSusp.Range("I$r") $donor50 26
Set-CellFromDonor $sheetSusp.Range("J$r") $donor50 69
Set-CellFromDonor $sheetSusp.Range("K$r") $donor50 187
Set-CellFromDonor $sheetSusp.Range("L$r") $donor14 "*"
Set-CellFromDonor $sheetSusp.Range("M$r") $donor14 "*"
Set-CellFromDonor $sheetSusp.Range("N$r") $donor50 23
Set-CellFromDonor $sheetSusp.Range("O$r") $donor14 "*"
Set-CellFromDonor $sheetSusp.Range("P$r") $donor37 "*"
Set-CellFromDonor $sheetSusp.Range("Q$r") $donor13 484

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Move each table's "last active cell" selection onto the new row, to
# mirror where the editor's cursor ended up after typing the update.
# ---------------------------------------------------------------------
$sheetCases.Range("A82").Select()
$sheetICU.Range("A71").Select()
$sheetSusp.Range("A63").Select()

# Leave the workbook looking at the Contents sheet, as it did before.
$sheetContents.Select()
$sheetContents.Range("A1").Select()
